# Update posts.xlsx after post
# The post at row 237 ("「卒業おめでとう」...") was removed from the source data,
# so delete that entire row; Excel will automatically shift all the
# subsequent rows (238..272) up by one and adjust the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(237).Delete()
